$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (E1) so the new
# header cells F1:H1 match the bold/bordered/centered style used by the
# other headers (style index 1 in the original workbook).
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New header labels
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# New boolean outlier-flag columns (rows 2-12)
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false

$ws.Range("F3").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = $false

$ws.Range("F4").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = $true

$ws.Range("F5").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = $false

$ws.Range("F6").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = $false

$ws.Range("F7").Value = $false
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = $false

$ws.Range("F8").Value = $false
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = $false

$ws.Range("F9").Value = $false
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = $true

$ws.Range("F10").Value = $false
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = $false

$ws.Range("F11").Value = $false
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = $false

$ws.Range("F12").Value = $false
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = $false

Write-Output "done"
